# RPA datasets push 2023-11-08
# The "LS머트리얼즈" listing (previously reported in row 9 of the
# "02_38커뮤니케이션(최근일자기준)" sheet) moved up to row 3 in the refreshed
# scrape, pushing the rows that used to sit between it and the top of the
# table down by one. Re-apply the refreshed snapshot for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# A:회사명  B:수요예측일  C:희망공모가(원)  D:확정공모가  E:공모금액(백만)  F:주간사
$rows = @(
    @("블루엠텍",       "2023.11.22~11.28", "15,000~19,000", "-", 21000, "하나증권,키움증권"),
    @("LS머트리얼즈",    "2023.11.22~11.28", "4,400~5,500",   "-", 64350, "키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권"),
    @("삼성스팩9호",     "2023.11.20~11.21", "2,000~2,000",   "-", 20000, "삼성증권"),
    @("교보스팩15호",    "2023.11.20~11.21", "2,000~2,000",   "-", 7000,  "교보증권"),
    @("케이엔에스",      "2023.11.16~11.22", "19,000~22,000", "-", 14250, "신영증권"),
    @("NH스팩30호",      "2023.11.15~11.16", "2,000~2,000",   "-", 16000, "NH투자증권"),
    @("와이바이오로직스", "2023.11.10~11.16", "9,000~11,000", "-", 13500, "유안타증권"),
    @("에이텀",          "2023.11.09~11.15", "23,000~30,000","-", 14950, "하나증권")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
